# Update cryptos list with latest Price (column D) and Volume(1h) (column E)
# values scraped from coinranking.com.
#
# Column D values are stored as text (e.g. "57.60", "1.00", "43.738.85" which
# isn't even a valid number), so a leading apostrophe is used to force a text
# entry and avoid Excel auto-converting/reformatting them as numbers (which
# would silently drop significant trailing zeros or mis-parse multi-dot
# strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.838.66'
$ws.Range("E2").Value = '  +4.96%  '
$ws.Range("D3").Value = '''2.283.95'
$ws.Range("E3").Value = '  +2.35%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''231.35'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("D7").Value = '''63.16'
$ws.Range("E7").Value = '  +4.54%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.424'
$ws.Range("E9").Value = '  +4.76%  '
$ws.Range("D10").Value = '''0.0952'
$ws.Range("E10").Value = '  +5.33%  '
$ws.Range("D11").Value = '''57.60'
$ws.Range("E11").Value = '  -0.94%  '
$ws.Range("D12").Value = '''26.21'
$ws.Range("E12").Value = '  +13.98%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").Value = '''2.624.01'
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").Value = '''15.86'
$ws.Range("D16").Value = '''5.96'
$ws.Range("E16").Value = '  +6.22%  '
$ws.Range("D17").Value = '''0.813'
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").Value = '''2.280.70'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").Value = '''43.725.80'
$ws.Range("E19").Value = '  +4.82%  '
$ws.Range("D20").Value = '''0.0₃0958'
$ws.Range("E20").Value = '  +6.06%  '
$ws.Range("D21").Value = '''73.15'
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("D22").Value = '''6.18'
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D23").Value = '''252.35'
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("D24").Value = '''2.61'
$ws.Range("E24").Value = '  +9.94%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = '''2.47'
$ws.Range("E26").Value = '  +8.71%  '
$ws.Range("D27").Value = '''9.87'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("D28").Value = '''171.23'
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").Value = '''0.139'
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("D30").Value = '''20.58'
$ws.Range("E30").Value = '  +3.18%  '
$ws.Range("D31").Value = '''1.45'
$ws.Range("E31").Value = '  +2.17%  '
$ws.Range("D32").Value = '''2.73'
$ws.Range("E32").Value = '  +3.29%  '
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("D34").Value = '''0.0705'
$ws.Range("E34").Value = '  +8.80%  '
$ws.Range("D35").Value = '''5.11'
$ws.Range("E35").Value = '  +1.25%  '
$ws.Range("D36").Value = '''4.73'
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("D37").Value = '''6.58'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = '''3.70'
$ws.Range("E38").Value = '  +1.84%  '
$ws.Range("E39").Value = '  -1.06%  '
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").Value = '''11.05'
$ws.Range("E42").Value = '  +28.89%  '
$ws.Range("D43").Value = '''8.60'
$ws.Range("E43").Value = '  +0.93%  '
$ws.Range("D44").Value = '''4.63'
$ws.Range("E44").Value = '  +5.65%  '
$ws.Range("E45").Value = '  -6.89%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").Value = '''0.0967'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").Value = '''98.28'
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("D49").Value = '''17.13'
$ws.Range("E49").Value = '  +3.22%  '
$ws.Range("D50").Value = '''1.486.37'
$ws.Range("E50").Value = '  +1.20%  '
$ws.Range("D51").Value = '''2.32'
$ws.Range("E51").Value = '  +1.67%  '
